$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122

# ---------------------------------------------------------------
# HTML content blocks for the PublishingPageContent column (F)
# ---------------------------------------------------------------
$htmlAboutUs = '<p>Fifth replenish upon. Years divide. I us called thing dry waters he itself. Female very she&#39;d, whales stars the darkness i Beast fruit that of. Two <strong>behold</strong> moving.</p>
<p>Moved yielding upon you every us beginning place sea <em>creature</em> him good.</p>
<p>There hath, first lights air that <em>him</em> dominion midst very. Abundantly is dominion face place forth. Sea said grass it divided stars divide.</p>
<p>One. Creeping shall. Fly that <em>they&#39;re</em> stars divide air second moveth winged.</p>
<p>Darkness meat all. Fruit evening our <strong>don&#39;t</strong> third you he blessed. Gathered.</p>
<p>Living image fruit from all can&#39;t can&#39;t beginning given place earth shall lights.</p>'
$htmlOurHistory = '<p>Brought deep. Likeness kind won&#39;t fly kind doesn&#39;t seas sixth in greater won&#39;t land, seasons, midst it. Winged. Let.</p>
<p>Firmament <em>of</em> gathering greater fowl had that isn&#39;t forth which every second seas was <em>fourth</em> seed itself, let a them. Their.</p>
<p>Replenish given them man make god forth life.</p>
<p>There female <strong>tree</strong> that may saw. Forth dry that subdue dry second, seed place moved own called give fruit you seed shall called don&#39;t them moving.</p>
<p>Wherein saw brought, beast thing saw saw fowl one bring beast that.</p>
<p>Make which be air, replenish greater form morning years void don&#39;t. To <strong>from</strong> lesser.</p>'
$htmlOurCulture = '<p>Isn&#39;t bearing fifth gathered was third land unto us brought image also had subdue thing fruit.</p>
<p>Seas first. Made from creature, image seasons void sea <em>they&#39;re</em> second.</p>
<p>May bearing two don&#39;t make fowl under. Abundantly whales days grass thing brought in.</p>
<p>Fill bring creeping there. Created made set second multiply first one you&#39;re so seas multiply tree void.</p>
<p>Divided the unto stars isn&#39;t, replenish divided God. Man from above. Spirit spirit bearing dry. Life behold good.</p>
<p>Creature given, heaven, given divide, have kind appear together, can&#39;t. Tree void don&#39;t and there they&#39;re have Whales. <em>Two</em> <em>lights</em> all heaven his all.</p>'
$htmlOurPhilosophy = '<p>Meat life replenish. One.</p>
<p>Seasons forth The blessed Fifth. After.</p>
<p>Night heaven <em>third</em> morning appear.</p>
<p>Midst female deep two.</p>
<p>I, deep face, years beast. Can&#39;t air.</p>
<p><strong>Subdue</strong> given replenish lesser.</p>
<p>They&#39;re seas there. One. Unto the fruit.</p>
<p>Firmament. Moving you&#39;ll open, <strong>lights</strong> beginning. Won&#39;t, signs.</p>
<p>Winged doesn&#39;t behold you&#39;ll.</p>
<p>Morning wherein light. Winged which fowl it.</p>
<p>In meat kind kind creeping all.</p>
<p>Seas signs moveth divided brought.</p>
<p>Fish living i gathered, fruit wherein unto fill.</p>
<p>Place <strong>rule</strong> night beast lesser signs male.</p>
<p>Whose years forth place, whose was.</p>
<p>Male. Very called. Over <strong>in</strong> god fourth have fruit good hath whales.</p>'
$htmlNews = '<p>Is <em>have</em> waters. First it created <strong>their</strong> sea sea years, behold god. Called.</p>
<p>All air. Years fifth over. Days was beast had sixth behold evening don&#39;t. Stars every set. Itself yielding man together of called.</p>
<p>Darkness land make set morning above won&#39;t that.</p>
<p>Beginning together form male fruit moveth bring first green fourth all creeping. Greater every likeness have fruitful <em>blessed</em> every also.</p>
<p>Shall. Lesser given saying, light creature had likeness <em>so</em> herb void beginning. For very land Female given. Thing, place don&#39;t one.</p>
<p>They&#39;re life creature. Light upon made evening won&#39;t night so meat, waters firmament let fill. The.</p>'
$htmlContactUs = '<p>Unto have place his. After days let replenish, life created so The. For given.</p>
<p>Firmament to for moving Beast. Kind fruitful set from there and he had sixth. Can&#39;t great. Two every saw <em>fill</em> first gathered.</p>
<p>Green creepeth beginning sixth third forth without. Be moved, make together shall, beast hath good creeping blessed saying cattle. They&#39;re created won&#39;t.</p>
<p>Lights man thing fill fruitful. Isn&#39;t evening wherein firmament them over void given a greater to greater sixth darkness.</p>
<p>Forth. Form also fowl two and day created.</p>
<p>Is you&#39;re fill void deep may moved moving said moved The evening were replenish. Place.</p>'

# ---------------------------------------------------------------
# Row 2 (existing row) - only DynamiteNavigation (E) and
# PublishingPageContent (F) change; everything else is identical.
# ---------------------------------------------------------------
$ws.Range("E2").Value = "About us"

# Give F2 the new wrap-text style (based on the plain-text style used
# by column B) before writing the long HTML value into it.
$ws.Range("B2").Copy()
$ws.Range("F2").PasteSpecial($xlPasteFormats)
$ws.Range("F2").WrapText = $true
$ws.Range("F2").Value = $htmlAboutUs
$ws.Rows.Item(2).RowHeight = 13.5

# ---------------------------------------------------------------
# New rows 3-7 (About us / Our history / Our culture / 
# Our philosophy / News / Contact us navigation entries)
# ---------------------------------------------------------------
# --- Row 3 ---
$ws.Range("B2").Copy()
$ws.Range("B3:E3").PasteSpecial($xlPasteFormats)
$ws.Range("H3").PasteSpecial($xlPasteFormats)
$ws.Range("J3").PasteSpecial($xlPasteFormats)
$ws.Range("F3").PasteSpecial($xlPasteFormats)
$ws.Range("F3").WrapText = $true
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial($xlPasteFormats)
$ws.Range("I2").Copy()
$ws.Range("I3").PasteSpecial($xlPasteFormats)
$ws.Range("L2:W2").Copy()
$ws.Range("L3:W3").PasteSpecial($xlPasteFormats)
$ws.Range("Y2").Copy()
$ws.Range("Y3").PasteSpecial($xlPasteFormats)

$ws.Range("A3").Value2 = 2
$ws.Range("B3").Value = '2_.000'
$ws.Range("C3").Value = 'Content Item'
$ws.Range("D3").Value = 'Our history'
$ws.Range("E3").Value = 'Our history'
$ws.Range("F3").Value = $htmlOurHistory
$ws.Range("G3").Value2 = 41936
$ws.Range("H3").Value = 'Main Menu'
$ws.Range("I3").Value2 = 2
$ws.Range("J3").Value = 'Approved'
$ws.Rows.Item(3).RowHeight = 13.5

# --- Row 4 ---
$ws.Range("B2").Copy()
$ws.Range("B4:E4").PasteSpecial($xlPasteFormats)
$ws.Range("H4").PasteSpecial($xlPasteFormats)
$ws.Range("J4").PasteSpecial($xlPasteFormats)
$ws.Range("F4").PasteSpecial($xlPasteFormats)
$ws.Range("F4").WrapText = $true
$ws.Range("G2").Copy()
$ws.Range("G4").PasteSpecial($xlPasteFormats)
$ws.Range("I2").Copy()
$ws.Range("I4").PasteSpecial($xlPasteFormats)
$ws.Range("L2:W2").Copy()
$ws.Range("L4:W4").PasteSpecial($xlPasteFormats)
$ws.Range("Y2").Copy()
$ws.Range("Y4").PasteSpecial($xlPasteFormats)

$ws.Range("A4").Value2 = 3
$ws.Range("B4").Value = '3_.000'
$ws.Range("C4").Value = 'Content Item'
$ws.Range("D4").Value = 'Our culture'
$ws.Range("E4").Value = 'Our culture'
$ws.Range("F4").Value = $htmlOurCulture
$ws.Range("G4").Value2 = 41936
$ws.Range("H4").Value = 'Main Menu'
$ws.Range("I4").Value2 = 3
$ws.Range("J4").Value = 'Approved'
$ws.Rows.Item(4).RowHeight = 13.5

# --- Row 5 ---
$ws.Range("B2").Copy()
$ws.Range("B5:E5").PasteSpecial($xlPasteFormats)
$ws.Range("H5").PasteSpecial($xlPasteFormats)
$ws.Range("J5").PasteSpecial($xlPasteFormats)
$ws.Range("F5").PasteSpecial($xlPasteFormats)
$ws.Range("F5").WrapText = $true
$ws.Range("G2").Copy()
$ws.Range("G5").PasteSpecial($xlPasteFormats)
$ws.Range("I2").Copy()
$ws.Range("I5").PasteSpecial($xlPasteFormats)
$ws.Range("L2:W2").Copy()
$ws.Range("L5:W5").PasteSpecial($xlPasteFormats)
$ws.Range("Y2").Copy()
$ws.Range("Y5").PasteSpecial($xlPasteFormats)

$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value = '4_.000'
$ws.Range("C5").Value = 'Content Item'
$ws.Range("D5").Value = 'Our philosophy'
$ws.Range("E5").Value = 'Our philosophy'
$ws.Range("F5").Value = $htmlOurPhilosophy
$ws.Range("G5").Value2 = 41936
$ws.Range("H5").Value = 'Main Menu'
$ws.Range("I5").Value2 = 4
$ws.Range("J5").Value = 'Approved'
$ws.Rows.Item(5).RowHeight = 13.5

# --- Row 6 ---
$ws.Range("B2").Copy()
$ws.Range("B6:E6").PasteSpecial($xlPasteFormats)
$ws.Range("H6").PasteSpecial($xlPasteFormats)
$ws.Range("J6").PasteSpecial($xlPasteFormats)
$ws.Range("F6").PasteSpecial($xlPasteFormats)
$ws.Range("F6").WrapText = $true
$ws.Range("G2").Copy()
$ws.Range("G6").PasteSpecial($xlPasteFormats)
$ws.Range("I2").Copy()
$ws.Range("I6").PasteSpecial($xlPasteFormats)
$ws.Range("L2:W2").Copy()
$ws.Range("L6:W6").PasteSpecial($xlPasteFormats)
$ws.Range("Y2").Copy()
$ws.Range("Y6").PasteSpecial($xlPasteFormats)

$ws.Range("A6").Value2 = 5
$ws.Range("B6").Value = '5_.000'
$ws.Range("C6").Value = 'Content Item'
$ws.Range("D6").Value = 'News'
$ws.Range("E6").Value = 'News'
$ws.Range("F6").Value = $htmlNews
$ws.Range("G6").Value2 = 41936
$ws.Range("H6").Value = 'Main Menu'
$ws.Range("I6").Value2 = 5
$ws.Range("J6").Value = 'Approved'
$ws.Rows.Item(6).RowHeight = 13.5

# --- Row 7 ---
$ws.Range("B2").Copy()
$ws.Range("B7:E7").PasteSpecial($xlPasteFormats)
$ws.Range("H7").PasteSpecial($xlPasteFormats)
$ws.Range("J7").PasteSpecial($xlPasteFormats)
$ws.Range("F7").PasteSpecial($xlPasteFormats)
$ws.Range("F7").WrapText = $true
$ws.Range("G2").Copy()
$ws.Range("G7").PasteSpecial($xlPasteFormats)
$ws.Range("I2").Copy()
$ws.Range("I7").PasteSpecial($xlPasteFormats)
$ws.Range("L2:W2").Copy()
$ws.Range("L7:W7").PasteSpecial($xlPasteFormats)
$ws.Range("Y2").Copy()
$ws.Range("Y7").PasteSpecial($xlPasteFormats)

$ws.Range("A7").Value2 = 6
$ws.Range("B7").Value = '6_.000'
$ws.Range("C7").Value = 'Content Item'
$ws.Range("D7").Value = 'Contact us'
$ws.Range("E7").Value = 'Contact us'
$ws.Range("F7").Value = $htmlContactUs
$ws.Range("G7").Value2 = 41936
$ws.Range("H7").Value = 'Main Menu'
$ws.Range("I7").Value2 = 6
$ws.Range("J7").Value = 'Approved'
$ws.Rows.Item(7).RowHeight = 13.5

# ---------------------------------------------------------------
# Column widths (E and F grew to accommodate the new content)
# ---------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 13
$ws.Columns.Item(6).ColumnWidth = 99.66666666666667

# ---------------------------------------------------------------
# Restore the on-screen selection to where the author left it
# ---------------------------------------------------------------
$ws.Range("F21").Select()
